# Update "想去人数" (interested-count) figures for three events.
# These figures are duplicated across the "展览" sheet (sheet1) and the
# "全部类型" aggregate sheet (sheet4).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAll        = $wb.Worksheets.Item("全部类型")

# 展览 sheet ("展览")
$wsExhibition.Range("F4").Value = 3318   # 合肥·第九届环形宇宙动漫游戏嘉年华
$wsExhibition.Range("F6").Value = 29     # 合肥·风月引代号鸢同人only
$wsExhibition.Range("F7").Value = 157    # 合肥·心动恋章·冬日序国乙&代号鸢同人only

# 全部类型 sheet (aggregate of all event types)
$wsAll.Range("F8").Value  = 3318   # 合肥·第九届环形宇宙动漫游戏嘉年华
$wsAll.Range("F10").Value = 29     # 合肥·风月引代号鸢同人only
$wsAll.Range("F12").Value = 157    # 合肥·心动恋章·冬日序国乙&代号鸢同人only
